$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# "Diametre trou bobine" (coil hole diameter) input changed from 8 to 10.5 mm;
# every dependent formula in the comparison table (H:U) and the base column
# (E) recalculates automatically from this single input edit.
$ws.Range("E6").Value = 10.5

# The highlighted (yellow) 2-decimal formatting that used to mark the
# "selected" column (H12:U12, wire-diameter row) was cleared back to the
# plain 2-decimal numeric style, with no fill.
$ws.Range("H12:U12").ClearFormats()
$ws.Range("H12:U12").NumberFormat = "0.00"

# Cursor / active-cell selection left on M36 when the file was saved.
$ws.Range("M36").Select() | Out-Null
